$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown

# Row 26 currently holds customer phone "09876543" (stored as text, leading zero kept)
# with an empty birthday and 0 total_points.
#
# The edit duplicates that row: row 26 keeps the same birthday/points but the phone
# becomes the numeric value 9876543 (leading zero dropped), while a new row 27 is
# appended underneath that preserves the original text phone "09876543" together with
# the same blank birthday / 0 points.
#
# We build this by copying the existing row 26 and inserting it (shifting the original
# row down to 27), which keeps row 27 byte-for-byte identical to the original row 26 -
# including its blank "birthday" cell - and then just overwrite the phone number in the
# (new) row 26 with the numeric value.
$ws.Rows(26).Copy()
$ws.Rows(26).Insert($xlShiftDown)

$ws.Cells.Item(26, 1).Value = 9876543
